$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2368
$ws.Range("I32").Value = 2810.1667
$ws.Range("J32").Value = 2036.375
$ws.Range("K32").Value = 2810.1667
$ws.Range("L32").Value = 2036.375
$ws.Range("M32").Value = -2484.1667
$ws.Range("N32").Value = -2688.375

$ws.Range("H58").Value = 5271.6665
$ws.Range("I58").Value = 2907.5
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 8722.5
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -8572.5
$ws.Range("N58").Value = -30300

$ws.Range("H137").Value = 568268.1
$ws.Range("I137").Value = 1672.9032
$ws.Range("J137").Value = 830424.1
$ws.Range("K137").Value = 5018.7096
$ws.Range("L137").Value = 2491272.3
$ws.Range("M137").Value = -2468.7096
$ws.Range("N137").Value = -2496372.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 201.33333
$ws.Range("I10").Value = 201.33333
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 201.33333
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -61.33332999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 163.61111
$ws.Range("I7").Value = 146.66667
$ws.Range("J7").Value = 172.08333
$ws.Range("K7").Value = 146.66667
$ws.Range("L7").Value = 172.08333
$ws.Range("M7").Value = -33.66667000000001
$ws.Range("N7").Value = -398.08333

$ws.Range("H22").Value = 596.0769
$ws.Range("I22").Value = 468.625
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 468.625
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -118.625
$ws.Range("N22").Value = -1500

$ws.Range("H107").Value = 1486.3
$ws.Range("I107").Value = 640.2
$ws.Range("J107").Value = 2332.4
$ws.Range("K107").Value = 640.2
$ws.Range("L107").Value = 2332.4
$ws.Range("M107").Value = 1279.8
$ws.Range("N107").Value = -6172.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1953.9166
$ws.Range("I17").Value = 165.22223
$ws.Range("J17").Value = 7320
$ws.Range("K17").Value = 495.66669
$ws.Range("L17").Value = 21960
$ws.Range("M17").Value = -326.66669
$ws.Range("N17").Value = -22298

$ws.Range("H34").Value = 298
$ws.Range("I34").Value = 298
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 894
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -810
$ws.Range("N34").ClearContents()

$ws.Range("H39").Value = 4249.8335
$ws.Range("I39").Value = 800
$ws.Range("J39").Value = 4939.8
$ws.Range("K39").Value = 2400
$ws.Range("L39").Value = 14819.4
$ws.Range("M39").Value = -2106
$ws.Range("N39").Value = -15407.4

$ws.Range("H55").Value = 3973.077
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 4137.5
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 12412.5
$ws.Range("M55").Value = -5823
$ws.Range("N55").Value = -12766.5

$ws.Range("H68").Value = 810.4
$ws.Range("I68").Value = 658.76666
$ws.Range("J68").Value = 875.3857400000001
$ws.Range("K68").Value = 1976.29998
$ws.Range("L68").Value = 2626.15722
$ws.Range("M68").Value = -1165.29998
$ws.Range("N68").Value = -4248.15722

$ws.Range("H71").Value = 810.4
$ws.Range("I71").Value = 658.76666
$ws.Range("J71").Value = 875.3857400000001
$ws.Range("K71").Value = 5928.89994
$ws.Range("L71").Value = 7878.47166
$ws.Range("M71").Value = -1872.89994
$ws.Range("N71").Value = -15990.47166

$ws.Range("H131").Value = 903.83905
$ws.Range("I131").Value = 398.1
$ws.Range("J131").Value = 969.5195
$ws.Range("K131").Value = 1194.3
$ws.Range("L131").Value = 2908.5585
$ws.Range("M131").Value = 3845.7
$ws.Range("N131").Value = -12988.5585

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121.94444
$ws.Range("I2").Value = 9.714286
$ws.Range("J2").Value = 193.36363
$ws.Range("K2").Value = 9.714286
$ws.Range("L2").Value = 193.36363
$ws.Range("M2").Value = 103.285714
$ws.Range("N2").Value = -419.36363

$ws.Range("H135").Value = 37352.94
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 37352.94
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 37352.94
$ws.Range("N135").Value = -47492.94

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8109.091
$ws.Range("I40").Value = 7837.5
$ws.Range("J40").Value = 8833.333000000001
$ws.Range("K40").Value = 7837.5
$ws.Range("L40").Value = 8833.333000000001
$ws.Range("M40").Value = -7701.5
$ws.Range("N40").Value = -9105.333000000001

$ws.Range("H122").Value = 3808
$ws.Range("I122").Value = 3510
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 10530
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -8080

$ws.Range("H136").Value = 1526.1833
$ws.Range("I136").Value = 1568.7435
$ws.Range("J136").Value = 1447.1428
$ws.Range("K136").Value = 4706.2305
$ws.Range("L136").Value = 4341.428400000001
$ws.Range("M136").Value = -2156.2305
$ws.Range("N136").Value = -9441.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1601.6666
$ws.Range("I13").Value = 1601.6666
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1601.6666
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1461.6666

$ws.Range("H51").Value = 18000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 18000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 18000
$ws.Range("N51").Value = -19020
$ws.Range("M51").ClearContents()

$ws.Range("H62").Value = 9776.799999999999
$ws.Range("I62").Value = 8628
$ws.Range("J62").Value = 11500
$ws.Range("K62").Value = 8628
$ws.Range("L62").Value = 11500
$ws.Range("M62").Value = -8004
$ws.Range("N62").Value = -12748

$ws.Range("H65").Value = 9776.799999999999
$ws.Range("I65").Value = 8628
$ws.Range("J65").Value = 11500
$ws.Range("K65").Value = 43140
$ws.Range("L65").Value = 57500
$ws.Range("M65").Value = -40020
$ws.Range("N65").Value = -63740

$ws.Range("H81").Value = 2148.2856
$ws.Range("I81").Value = 1696
$ws.Range("J81").Value = 2487.5
$ws.Range("K81").Value = 3392
$ws.Range("L81").Value = 4975
$ws.Range("M81").Value = -2331
$ws.Range("N81").Value = -7097

$ws.Range("H84").Value = 2148.2856
$ws.Range("I84").Value = 1696
$ws.Range("J84").Value = 2487.5
$ws.Range("K84").Value = 16960
$ws.Range("L84").Value = 24875
$ws.Range("M84").Value = -11656
$ws.Range("N84").Value = -35483

$ws.Range("H119").Value = 30000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 30000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676

$ws.Range("H122").Value = 135745
$ws.Range("I122").Value = 2950
$ws.Range("J122").Value = 241981
$ws.Range("K122").Value = 8850
$ws.Range("L122").Value = 725943
$ws.Range("M122").Value = -6400
$ws.Range("N122").Value = -730843
